# FieldRenameMap.xlsx edit script
# Adds 16 new field-rename rows (rows 48-63) to Sheet1, extends the
# "%"+A formula down column B (as a shared formula for the B49:B63 block,
# matching how Excel fills a formula down a selection), sets column C to
# "Y" for each new row, and re-splits the existing conditional formatting
# rule on column C so the newly added rows keep the same "Y" highlight
# rule while rows 54 and 55 each pick up their own (identical) rule - this
# mirrors the dxf/rule split seen after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows 48-63 -----------------------------------------------
$names = @(
    "EX.Comment",
    "EX.Field",
    "IN.Comment",
    "IN.Comment",
    "TR.Comment",
    "TR.Target Field Definition",
    "IN.Target Field Order",
    "IN.Target Field Name",
    "Source Folder",
    "MD Field from Source",
    "MD Field Name in This Datamodel",
    "MD Field Name in This Datamodel Tmp",
    "MD QVD Read",
    "MD Source Folder",
    "MD Source To QVD",
    "MD Table Name in This Datamodel"
)

$startRow = 48
for ($i = 0; $i -lt $names.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $names[$i]
    $ws.Cells.Item($r, 3).Value = "Y"
}

# Row 48's B formula is entered on its own (not part of the fill-down).
$ws.Range("B48").Formula = '="%"&A48'

# Rows 49-63 are filled down together from B49, producing one shared
# formula block (B49:B63), same as dragging the fill handle down.
$ws.Range("B49:B63").Formula = '="%"&A49'

# --- Conditional formatting on column C ---------------------------------
# Re-home the original "Y" rule (still dxf #0) onto just the C56:C63
# block, then add matching new rules (same green fill / dark green font)
# for C54, C55 and the remaining C1:C53 region, in that order, so the
# priorities / dxf indices come out 1..4 as Excel would number them.
$greenFont = 24832      # RGB(0,97,0)   -> FF006100
$greenFill = 13561798   # RGB(198,239,206) -> FFC6EFCE

$originalRule = $ws.Range("C1:C1048576").FormatConditions.Item(1)
$originalRule.ModifyAppliesToRange($ws.Range("C56:C63"))

$ruleC55 = $ws.Range("C55").FormatConditions.Add(1, 3, '"Y"')
$ruleC55.Font.Color = $greenFont
$ruleC55.Interior.Color = $greenFill

$ruleC54 = $ws.Range("C54").FormatConditions.Add(1, 3, '"Y"')
$ruleC54.Font.Color = $greenFont
$ruleC54.Interior.Color = $greenFill

$ruleRest = $ws.Range("C1:C53").FormatConditions.Add(1, 3, '"Y"')
$ruleRest.Font.Color = $greenFont
$ruleRest.Interior.Color = $greenFill

# --- View state -----------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$ws.Range("D55").Select() | Out-Null
